$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# --- Moisturizer test-case block --------------------------------------
# The old block had three data rows: Y/Chrome/Aloe/Least, Y/Chrome/Almond/Least
# and N/Mozilla/Aloe/Expensive. The new testcase selects a moisturizer and
# adds it to the cart by content (Aloe) AND content2 (Almond) in one row, so
# drop the extra "Almond" row (was row 8) and the "N/Mozilla" row (now
# shifted up to row 8).
$ws.Rows(8).Delete()
$ws.Rows(8).Delete()

# --- Sunscreen test-case block ------------------------------------------
# Likewise drop its "N/Mozilla" row (now shifted up to row 12).
$ws.Rows(12).Delete()

# Moisturizer header/data: split the single "ProductContent" column into two
# columns (ProductContent1 / ProductContent2) and record both product
# contents (Aloe, Almond) being added to the cart.
$ws.Range("C6").Value = "ProductContent1"
$ws.Range("C10").Value = "ProductContent1"
$ws.Range("D6").Value = "ProductContent2"
$ws.Range("D10").Value = "ProductContent2"
$ws.Range("D7").Value = "Almond"

# Sunscreen header/data: same ProductContent1/ProductContent2 split, now
# recording two SPF options (SPF-50, SPF-30) added to the cart.
$ws.Range("C11").Value = "SPF-50"
$ws.Range("D11").Value = "SPF-30"

# Widen columns C:D to fit the new, longer header text.
$ws.Columns("C:D").AutoFit() | Out-Null

# Leave the selection where the last edit happened, matching the workbook's
# recorded cursor position.
$ws.Range("D11").Select() | Out-Null
